$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the original data block (A1:C34, which includes formulas in column A)
# and paste VALUES ONLY (xlPasteValues) into the new location D9:F42 -
# this is how the shared-formula column (A3:A34) turns into plain numbers.
$src = $ws.Range("A1:C34")
$src.Copy()
$dst = $ws.Range("D9")
$dst.PasteSpecial(-4163)

$excel.CutCopyMode = 0

# The data now lives at D9:F42, so clear out the old A1:C34 block.
$ws.Range("A1:C34").Clear()

# Update the selection to match the new active cell.
$ws.Range("F9").Select()
